# Regenerate s_vals data to filter save games: update B:G values for rows 2-8.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.445647641019636, 0.04103571897497393, 0.1496068669990043, 0.5333859586016987, 2.169676185595313)
    3 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    4 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    5 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    6 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    7 = @(0.2881169905109251, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.598097515653722)
    8 = @(0.1169995834814548, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 1.67637130870356)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G
}
